# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") values for rows 2-33 on the active sheet,
# replacing the previously scraped "Strike#" derived figures with the
# recomputed "K" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 1
    6  = 5
    7  = 0
    8  = 1
    9  = 2
    10 = 2
    11 = 3
    12 = 0
    13 = 1
    14 = 1
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 3
    20 = 1
    21 = 2
    22 = 5
    23 = 3
    24 = 3
    25 = 1
    26 = 1
    27 = 3
    28 = 0
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
